# Adding the number of differences ("# of Diffs") as a new first column
# to the (normal layout) report sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts the existing
# ID / FIRST_NAME (Source1) / FIRST_NAME (Source2) columns one to the right.
$ws.Columns("A:A").Insert()

# The old "ID" column (now column B) used to carry the red/underlined
# style used for the data rows; that styling now belongs to the new
# "# of Diffs" column instead, so copy it over and reset column B to
# the default (unstyled) look.
$ws.Range("C2:C7").Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)
$ws.Range("B2:B7").Style = "Normal"

# Header row.
$ws.Range("A1").Value = "# of Diffs"

# Data rows: every row currently has exactly one diff.
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = 1.0
}
